# Edit script: expand paragraph 1 with the first part of the recording's
# continuation, correct the transcription text in paragraph 2, split off a
# new paragraph for "lo siguiente: ..." and leave a trailing empty paragraph.

$d = $word.ActiveDocument

# --- Step 1: extend paragraph 1 (merge in corrected continuation text) ---
$old1 = "la ansiedad del siglo XXI es más conocida como el estrés. y Dios no nos manda a vivir con estrés, sino nos manda a vivir en paz. y la paz no es algo externo, sino que está dentro de nuestro corazón y viene a consecuencia de confiar nosotros en Dios. el mundo produce paz. una puesta de sol puede producirlos. Paz también lo puede hacer. el regar tus plantitas te puede producir paz, una sensación de calma, el jugar con tu mascota en sacarla a pasear, el echarnos en una hamaca, a mí me produce mucha paz. pero podemos estar lo más relajados del mundo, y viene alguien, nos dan una mala noticia y nos caemos de la manta, dejamos de regar las plantas y se acabó todo."
$new1 = "la ansiedad del siglo XXI es más conocida como el estrés. y Dios no nos manda a vivir con estrés, sino nos manda a vivir en paz. y la paz no es algo externo, sino que está dentro de nuestro corazón y viene a consecuencia de confiar nosotros en Dios. el mundo produce paz. una puesta de sol puede producirlos. Paz también lo puede hacer. el regar tus plantitas te puede producir paz, una sensación de calma, el jugar con tu mascota en sacarla a pasear, el echarnos en una hamaca, a mí me produce mucha paz. pero podemos estar lo más relajados del mundo y viene alguien, nos dan una mala noticia y nos caemos de la manta, dejamos de regar las plantas y se acabó todo. porque la paz que el mundo nos da una paz temporal, mientras que la paz de Dios nos da en cambio paz a pesar de las circunstancias, a pesar de lo que dice la gente. la gente dirá como ella, a pesar de estar pasando por problemas, sigue teniendo paz, puede seguir reflejando paz. como ella, a pesar de de de las necesidades económicas que tiene, no está quejándose ni está renegando. si no está compás, cómo es que ya que, a pesar de sus dolores, se da el tiempo para preocuparse y atender de otros y dice la Biblia va hasta el de. a cada día su propio pan. y si nos ponemos a mirar, la mayor parte de nuestras preocupaciones no son acerca del hoy, sino son acerca del mañana. nos estamos preocupando por lo que sucederá mañana, por cosas que nosotros no tenemos en nuestro control. la Biblia dice que no debemos estar afanosos por nada sino que sean conocidas nuestras peticiones en toda oración y suplica en el espíritu. así que todos los días nosotros necesitamos pedirle paz a Dios. cuando las preocupaciones legítimas de nuestra vida son manejadas erróneamente, pueden llevarnos a tener preocupaciones dominantes que nos conduzcan al temor. por eso la Biblia nos habla de esto en job, capítulo 3, verso 25: lo que temí me aconteció. para el que cree, en realidad, todo es posible. si tú piensas: no, yo no voy a poder, no vas a poder."
$r1 = $d.Content
$found1 = $r1.Find.Execute($old1, $false, $false, $false, $false, $false, $true, 1, $false, $new1, 2)
if (-not $found1) {
    throw "Step 1 replacement (paragraph 1 text) was not found"
}

# --- Step 2: rewrite paragraph 2's text, split it into two paragraphs, and
#     leave a new empty paragraph at the end ---
$old2 = " porque la paz que el mundo nos da una paz temporal, mientras que la paz de Dios nos da en cambio paz a pesar de las circunstancias, a pesar de lo que dice la gente. la gente dirá como ella, a pesar de estar pasando por problemas, sigue teniendo paz. puedes seguir reflejando paz como ella, a pesar de de de las necesidades económicas que tiene, no está quejándose ni está renegando. si no está, compás, cómo es que ya que, a pesar de sus dolores, se da el tiempo para preocuparse y atender de otros y dice: la Biblia va hasta el de a cada día su propio pan? y si nos ponemos a mirar, la mayor parte de nuestras preocupaciones no son acerca de hoy, sino son acerca del mañana. nos estamos preocupando por lo que sucederá mañana, por cosas que nosotros no tenemos en nuestro control. la Biblia dice que no debemos estar afanosos por nada, sino que sean conocidas nuestras peticiones en toda oración y suplica en el espíritu, así que todos los días nosotros necesitamos. pedirle paz a Dios. cuando las preocupaciones legítimas de nuestra vida son manejadas erróneamente, pueden llevarnos a tener preocupaciones dominantes que nos conduzcan al temor. por eso la Biblia nos habla de esto en job, capítulo 3, verso 25. job decía: lo que temí aconteció. para el que cree en realidad, todo es posible. si tú piensas: no, yo no voy a poder, no vas a poder, no creo que este mes me alcance el dinero, entonces no te va a alcanzar porque nosotros tenemos en nuestra lengua el poder para dar vida o muerte, tenemos el poder de las palabras, de lo que nosotros decimos. apenas nos despertamos, le pedimos a Dios paz y decimos que, como nuestros días serán nuestras fuerzas, que Dios suplirá toda necesidad, que ese problema no nos alterarnos nervios, porque nuestros nervios están en las manos de Dios. entonces, qué sucede cuando nosotros nos dejamos abrumar y comenzamos y permitimos que el estrés tome control de nuestras vidas y no la paz de Dios? cómo es que nos captura entonces la preocupación cuando comenzamos a dibujar escenarios en nuestra mente que el futuro comienzan a ahogar nuestro presente y comenzamos a desconfiar de las promesas de Dios, del poder de Dios? el libro de Lucas, capítulo 8, verso 14, dice la palabra, porque acá nos están hablando de la semilla de la Palabra? dice la que cayó entre espinos. representan los que oyen, pero después de un tiempo los ahogan las preocupaciones, las riquezas, los placeres de la vida y no llegan a madurar. entonces nosotros podemos haber recibido una promesa, podemos haber recibido una palabra, pero con el tiempo esa palabra que entró por este oído sale por el otro y comienza a ahogarse con las preocupaciones. las riquezas y los seres de la vida y no recibimos la promesa que Dios tiene para nosotros. entonces qué daño nos causa? lo primero es que roba, roba nuestra concentración y nos hace improductivos. porque cuando estamos nosotros estresados, no podemos concentrarnos en una sola cosa, no avanzamos, no producimos, nos seguimos adelante, sino que estamos concentrados en esa sola cosa, en esa sola cosa, en esa sola cosa. le damos vueltas solamente a esa cosa. entonces roba nuestra concentración y nos hace improductivos. lo siguiente: roba nuestro descanso. por qué? porque podríamos estar descansando en ese momento, pero el estrés es como un "
$new2 = " no creo que este mes me alcance el dinero, entonces no te va a alcanzar, porque nosotros tenemos, en nuestra lengua, el poder para dar vida o muerte. tenemos el poder de las palabras, de lo que nosotros decimos. apenas nos despertamos, le pedimos a Dios paz y decimos que, como nuestros días serán nuestras fuerzas, que Dios suplirá toda necesidad, que ese problema no nos va a alterar los nervios, porque nuestros nervios están en las manos de Dios. entonces, qué sucede cuando nosotros nos dejamos abrumar y comenzamos y permitimos que el estrés tome control de nuestras vidas y no la paz de Dios? cómo es que nos captura entonces la preocupación cuando comenzamos a dibujar escenarios en nuestra mente que el futuro comienzan a ahogar nuestro presente y comenzamos a desconfiar de las promesas de Dios, del poder de Dios? el libro de Lucas, capítulo 8, verso 14, dice la palabra, porque acá nos están hablando de la semilla de la Palabra. dice la que cayó entre espinos, representa a los que oyen. pero después de un tiempo los ahogan las preocupaciones, las riquezas, los placeres de la vida y no llegan a madurar. entonces nosotros podemos haber recibido una promesa, podemos haber recibido una palabra, pero con el tiempo esa palabra que entró por este oído sale por el otro y comienza a ahogarse con las preocupaciones, las riquezas y los placeres de la vida y no recibimos la promesa que Dios tiene para nosotros. entonces qué daño nos causa? lo primero es que roba, roba nuestra concentración y nos hace improductivos. porque cuando estamos nosotros estresados, no podemos concentrarnos en una sola cosa, no avanzamos, no producimos, nos seguimos adelante. sino que estamos concentrados en esa sola cosa, en esa sola cosa, en esa sola cosa le damos vueltas solamente a esa cosa, entonces roba nuestra concentración y nos hace improductivos.^p lo siguiente: robo a nuestro descanso. por qué? porque podríamos estar descansando en ese momento, pero el estrés es como un^p"
$r2 = $d.Content
$found2 = $r2.Find.Execute($old2, $false, $false, $false, $false, $false, $true, 1, $false, $new2, 2)
if (-not $found2) {
    throw "Step 2 replacement (paragraph 2 text) was not found"
}

Write-Output "Paragraphs: $($d.Paragraphs.Count)"
